$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto market data (prices + 1h volume %) scraped on
# Fri Mar 31 20:57:19 UTC 2023. Rows 10/11, 31/32 and 34/35 also swap
# rank order (coin name + link) in addition to the value refresh.
# Price cells that look like plain decimals are written with a leading
# apostrophe so Excel keeps storing them as text (matching the sheet's
# existing text-formatted price column) instead of coercing to a number.

$ws.Range("D2").Value = '28.460.88'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").Value = '1.824.72'
$ws.Range("E3").Value = '  +1.68%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'316.23"
$ws.Range("E5").Value = '  -0.19%  '

$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").Value = "'0.5467"
$ws.Range("E7").Value = '  +2.33%  '

$ws.Range("D8").Value = "'0.4037"
$ws.Range("E8").Value = '  +7.17%  '

$ws.Range("D9").Value = "'0.07692"
$ws.Range("E9").Value = '  +3.02%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = "'1.116"
$ws.Range("E10").Value = '  +1.85%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = "'41.89"
$ws.Range("E11").Value = '  +0.16%  '

$ws.Range("E12").Value = '  +3.52%  '

$ws.Range("D13").Value = "'7.664"
$ws.Range("E13").Value = '  +5.93%  '

$ws.Range("D14").Value = "'1.001"
$ws.Range("E14").Value = '  +0.09%  '

$ws.Range("D15").Value = "'20.92"
$ws.Range("E15").Value = '  +1.30%  '

$ws.Range("D16").Value = '1.828.68'
$ws.Range("E16").Value = '  +2.38%  '

$ws.Range("D17").Value = "'0.00001082"
$ws.Range("E17").Value = '  +2.41%  '

$ws.Range("D18").Value = "'89.77"
$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("D19").Value = "'0.06597"
$ws.Range("E19").Value = '  +2.03%  '

$ws.Range("E20").Value = '  +1.75%  '

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").Value = "'6.064"
$ws.Range("E22").Value = '  +2.60%  '

$ws.Range("D23").Value = '28.467.74'
$ws.Range("E23").Value = '  +1.17%  '

$ws.Range("D24").Value = "'11.11"
$ws.Range("E24").Value = '  -0.86%  '

$ws.Range("E25").Value = '  +6.54%  '

$ws.Range("D26").Value = "'2.473"
$ws.Range("E26").Value = '  +7.86%  '

$ws.Range("D27").Value = "'20.75"
$ws.Range("E27").Value = '  +2.48%  '

$ws.Range("D28").Value = "'157.39"
$ws.Range("E28").Value = '  +1.68%  '

$ws.Range("D29").Value = '2.039.47'
$ws.Range("E29").Value = '  +2.47%  '

$ws.Range("D30").Value = "'123.94"
$ws.Range("E30").Value = '  +2.92%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'1.127"
$ws.Range("E31").Value = '  +1.02%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = "'0.1109"
$ws.Range("E32").Value = '  +6.07%  '

$ws.Range("D33").Value = "'5.679"
$ws.Range("E33").Value = '  +1.99%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = "'3.643"
$ws.Range("E34").Value = '  -0.23%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.07335"
$ws.Range("E35").Value = '  +12.58%  '

$ws.Range("D36").Value = "'0.2255"
$ws.Range("E36").Value = '  -0.38%  '

$ws.Range("D37").Value = "'0.02341"
$ws.Range("E37").Value = '  +2.40%  '

$ws.Range("D38").Value = "'5.206"
$ws.Range("E38").Value = '  +3.64%  '

$ws.Range("D39").Value = "'8.864"
$ws.Range("E39").Value = '  +4.49%  '

$ws.Range("D40").Value = "'11.37"
$ws.Range("E40").Value = '  +2.63%  '

$ws.Range("D41").Value = "'0.6276"
$ws.Range("E41").Value = '  +1.85%  '

$ws.Range("D42").Value = "'1.180"
$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("D44").Value = "'1.394"
$ws.Range("E44").Value = '  -3.77%  '

$ws.Range("D45").Value = "'13.50"
$ws.Range("E45").Value = '  +1.72%  '

$ws.Range("D46").Value = "'3.698"
$ws.Range("E46").Value = '  +0.70%  '

$ws.Range("D47").Value = "'0.5851"
$ws.Range("E47").Value = '  +1.33%  '

$ws.Range("D48").Value = "'125.02"
$ws.Range("E48").Value = '  -1.78%  '

$ws.Range("E49").Value = '  +3.88%  '

$ws.Range("D50").Value = "'1.202"
$ws.Range("E50").Value = '  +1.18%  '

$ws.Range("D51").Value = "'0.06884"
$ws.Range("E51").Value = '  +1.05%  '

